$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the reduction index values in column A (rows 2-8) from 1..7 to 0..6
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6

# Update the selected cell/range on the sheet
$ws.Range("B13").Select()
